# Applies the "Estado de Cuenta" update described by the diff:
#  - Updates the total "Valor Mora" amount and period count.
#  - The existing period rows shift: what used to be the newest period (2505)
#    is replaced by a new newest period (2507); the old middle period (2506)
#    keeps its value but becomes a normal row; a new last row is added for
#    the old oldest period (2505) which keeps its original amount.
#  - Inserts a new row (old last row to become a normal row, and a brand
#    new last row is appended) before the signature block, which shifts the
#    signature line / label rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update total "Valor Mora" (E11) and "Cant. Periodos" count (F13)
$ws.Range("E11").Value = 168630
$ws.Range("F13").Value = 3

# 2) Insert a new row at 18 (pushes the signature rows 22/23 down to 23/24,
#    and shifts the merged cell ranges together with them).
$ws.Rows(18).Insert()

# 3) The freshly inserted row 18 has no formatting yet - give it the same
#    look the old "last" detail row (17) had, then fill it with the data
#    for the oldest period (2505), which keeps its original amount.
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1002202627"
$ws.Range("D18").Value = "EDUARDO RAFAEL FIGUEROA ROMERO"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 37230
$ws.Range("G18").Value = 1642500

# 4) Row 17 (previously the last detail row, holding period 2506) now sits
#    in the middle of the table, so give it the "middle row" formatting
#    that row 16 uses. Its data (2506 / 65700) stays the same.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# 5) Row 16 becomes the newest period (2507) with its new amount.
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 65700

$excel.CutCopyMode = 0
